# Migrate mountebank stub example to use api template framework.
# Add a new "Create Stub" template row (mountebank / mb) to the Sheet1 template list,
# update column sizing / zoom / selection to match the author's final view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row of template data (row 5) ---------------------------------
# (cell order chosen to match the shared-string allocation order of the target file)
$ws.Range("C5").Value = "mb"
$ws.Range("D5").Value = "Create Stub"
$ws.Range("F5").Value = "test"
$ws.Range("B5").Value = '{"optionalField": true}'
$ws.Range("A5").Value = 4
$ws.Range("E5").Value = "JSON"

# Match the wrap-text style already used for the "Template"/"ProjectName" columns
# on the other data rows (cellXfs index 1 = wrapText).
$ws.Range("B5:C5").WrapText = $true

# --- Re-size rows / columns to match the refreshed template layout -------------
$ws.Rows.Item(2).RowHeight = 409.6
$ws.Rows.Item(3).RowHeight = 409.6
$ws.Rows.Item(4).RowHeight = 409.6
$ws.Rows.Item(5).RowHeight = 16

$ws.Columns.Item(1).ColumnWidth = 9.5
$ws.Columns.Item(2).ColumnWidth = 70.83
$ws.Columns.Item(3).ColumnWidth = 10.67
$ws.Columns.Item(4).ColumnWidth = 17.5
$ws.Columns.Item(6).ColumnWidth = 32.83

# --- Update the view: zoom in and select the next empty cell -------------------
$ws.Range("B6").Select()
$excel.ActiveWindow.Zoom = 125
